# This script applies a data update that corresponds to swapping the contents
# of several row pairs in the "Artfynd" worksheet (rows 5<->6, 9<->10,
# 13<->14, 18<->19). Only the cells whose values actually differ between the
# two rows of each pair are touched; the remaining cells in those rows are
# identical between the paired rows already, so they do not need to change.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row pair 5 <-> 6 ---
$ws.Range("A5").Value = 130806161
$ws.Range("Q5").Value = 472150
$ws.Range("R5").Value = 6864539

$ws.Range("A6").Value = 130806143
$ws.Range("Q6").Value = 472165
$ws.Range("R6").Value = 6864552

# --- Row pair 9 <-> 10 ---
$ws.Range("A9").Value = 130806878
$ws.Range("B9").Value = 80348
$ws.Range("E9").Value = 6458
$ws.Range("F9").Value = "Lunglav"
$ws.Range("G9").Value = "Lobaria pulmonaria"
$ws.Range("H9").Value = "(L.) Hoffm."
$ws.Range("Q9").Value = 472144
$ws.Range("R9").Value = 6864413
$ws.Range("S9").Value = 20

$ws.Range("A10").Value = 130806196
$ws.Range("B10").Value = 79243
$ws.Range("E10").Value = 6425
$ws.Range("F10").Value = "Garnlav"
$ws.Range("G10").Value = "Alectoria sarmentosa"
$ws.Range("H10").Value = "(Ach.) Ach."
$ws.Range("Q10").Value = 472131
$ws.Range("R10").Value = 6864526
$ws.Range("S10").Value = 10

# --- Row pair 13 <-> 14 ---
$ws.Range("A13").Value = 130839368
$ws.Range("B13").Value = 80349
$ws.Range("E13").Value = 2081
$ws.Range("F13").Value = "Skrovellav"
$ws.Range("G13").Value = "Lobaria scrobiculata"
$ws.Range("H13").Value = "(Scop.) DC."
$ws.Range("P13").Value = "Kristinehamnskojan, Kristinehamnskojan, Hjd"
$ws.Range("Q13").Value = 472174
$ws.Range("R13").Value = 6864517
$ws.Range("S13").Value = 10

$ws.Range("A14").Value = 130839206
$ws.Range("B14").Value = 80348
$ws.Range("E14").Value = 6458
$ws.Range("F14").Value = "Lunglav"
$ws.Range("G14").Value = "Lobaria pulmonaria"
$ws.Range("H14").Value = "(L.) Hoffm."
$ws.Range("P14").Value = "Öratjärnknallarna, Öratjärnknallarna, Hjd"
$ws.Range("Q14").Value = 472147
$ws.Range("R14").Value = 6864552
$ws.Range("S14").Value = 20

# --- Row pair 18 <-> 19 ---
$ws.Range("A18").Value = 130839361
$ws.Range("P18").Value = "Kristinehamnskojan, Kristinehamnskojan, Hjd"
$ws.Range("Q18").Value = 472176
$ws.Range("R18").Value = 6864518
$ws.Range("S18").Value = 10

$ws.Range("A19").Value = 130839096
$ws.Range("P19").Value = "Öratjärnknallarna, Öratjärnknallarna, Hjd"
$ws.Range("Q19").Value = 472178
$ws.Range("R19").Value = 6864559
$ws.Range("S19").Value = 20
